# Fruta / hortaliza, semanal
# Insert 2 new weekly records (rows) right before current row 478,
# pushing all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of this block.
$ws.Rows("478:479").Insert()

# --- New row 478 ---
$ws.Cells.Item(478, 1).Value  = 7
$ws.Cells.Item(478, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(478, 3).Value  = "Ñuble"
$ws.Cells.Item(478, 4).Value  = 45211
$ws.Cells.Item(478, 5).Value  = 16
$ws.Cells.Item(478, 6).Value  = "Fruta"
$ws.Cells.Item(478, 7).Value  = 100101
$ws.Cells.Item(478, 8).Value  = "Berries"
$ws.Cells.Item(478, 9).Value  = 100112025
$ws.Cells.Item(478, 10).Value = "Frutilla"
$ws.Cells.Item(478, 11).Value = "Sin especificar"
$ws.Cells.Item(478, 12).Value = "Primera"
$ws.Cells.Item(478, 13).Value = 150
$ws.Cells.Item(478, 14).Value = 11000
$ws.Cells.Item(478, 15).Value = 11000
$ws.Cells.Item(478, 16).Value = 11000
$ws.Cells.Item(478, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(478, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(478, 19).Value = 1571
$ws.Cells.Item(478, 20).Value = 7

# --- New row 479 ---
$ws.Cells.Item(479, 1).Value  = 7
$ws.Cells.Item(479, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(479, 3).Value  = "Ñuble"
$ws.Cells.Item(479, 4).Value  = 45211
$ws.Cells.Item(479, 5).Value  = 16
$ws.Cells.Item(479, 6).Value  = "Fruta"
$ws.Cells.Item(479, 7).Value  = 100101
$ws.Cells.Item(479, 8).Value  = "Berries"
$ws.Cells.Item(479, 9).Value  = 100112025
$ws.Cells.Item(479, 10).Value = "Frutilla"
$ws.Cells.Item(479, 11).Value = "Sin especificar"
$ws.Cells.Item(479, 12).Value = "Segunda"
$ws.Cells.Item(479, 13).Value = 150
$ws.Cells.Item(479, 14).Value = 8000
$ws.Cells.Item(479, 15).Value = 8000
$ws.Cells.Item(479, 16).Value = 8000
$ws.Cells.Item(479, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(479, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(479, 19).Value = 1143
$ws.Cells.Item(479, 20).Value = 7
